# Formed the consolidated report
# Fill in the "Absent" (column H) values that were previously left blank
# or incorrect, completing the consolidated attendance report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 1
$ws.Range("H16").Value = 0
